# Record the payments made by Aniket, Athul, Riyas, Tintu and Vikram.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Aniket (row 5) paid for June,18 (column F)
$ws.Range("F5").Value = 300

# Athul (row 6) paid for July,18 (column G)
$ws.Range("G6").Value = 500

# Riyas (row 20) paid for July,18 (column G)
$ws.Range("G20").Value = 500

# Tintu (row 27) paid for July,18 (column G)
$ws.Range("G27").Value = 500

# Vikram (row 31) had 400 recorded for July,18 (column G); correct it to 500
$ws.Range("G31").Value = 500

# Leave the selection where the user last clicked while entering data
$ws.Range("G20").Select()
